$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.137.31'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.140.38'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.56'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.55'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.137.90'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.87'
$ws.Range('E11').Value = '  +2.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  -3.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.20'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.121'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.653.49'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.879.74'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.137.43'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.13'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.55'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.32'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.731'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.41'
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.33'
$ws.Range('E24').Value = '  +7.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.96'
$ws.Range('E25').Value = '  -2.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '80.78'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.62'
$ws.Range('E28').Value = '  +10.78%  '
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.20'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.16'
$ws.Range('E32').Value = '  +4.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.09'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0857'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('E37').Value = '  -3.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.05'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.26'
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.53'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '440.79'
$ws.Range('E41').Value = '  -3.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.96'
$ws.Range('E42').Value = '  +2.65%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.286'
$ws.Range('E43').Value = '  +2.95%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0371'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.904.95'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.91'
$ws.Range('E46').Value = '  +15.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.108'
$ws.Range('E47').Value = '  -3.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.43'
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.21'
$ws.Range('E51').Value = '  +0.67%  '
